# Commit: "Hide 9650 as a service from service portfolio on devsite"
#
# The "Return mailbox parcel" service (service/request code 9650) is
# removed entirely from the Booking & SG API table. Deleting its row
# shifts every subsequent row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Cells.Find("Return mailbox parcel")
if ($target) {
    $target.EntireRow.Delete()
}

# Re-anchor the table's AutoFilter range now that the table is one row
# shorter (was A1:P66, now A1:P65).
$ws.AutoFilterMode = $false
$ws.Range("A1:P65").AutoFilter() | Out-Null

# The AutoFilter's backing defined name (_xlnm._FilterDatabase) keeps the
# stale upper bound after EntireRow.Delete/AutoFilter re-apply, so update
# it explicitly to match the new table extent.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Booking & SG API'!`$A`$1:`$P`$65"
    }
}
